# Auto-generated edit script: updates cryptos list (prices, volumes, and two
# row re-orderings) to match the Mon Oct 28 10:46:00 UTC 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''68.751.74'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.49%  '

# Row 3
$ws.Range("D3").Value = '''2.534.97'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.74%  '

# Row 4
$ws.Range("D4").Value = '''1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").Value = '''594.21'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.16%  '

# Row 6
$ws.Range("D6").Value = '''177.85'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.22%  '

# Row 7
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
$ws.Range("E8").Value = '  +1.33%  '

# Row 9
$ws.Range("D9").Value = '''2.534.22'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.71%  '

# Row 10
$ws.Range("E10").Value = '  +6.32%  '

# Row 11
$ws.Range("E11").Value = '  -0.99%  '

# Row 12
$ws.Range("E12").Value = '  +1.28%  '

# Row 13
$ws.Range("D13").Value = '''0.340'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.11%  '

# Row 15
$ws.Range("D15").Value = '''26.15'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.20%  '

# Row 16
$ws.Range("D16").Value = '''68.464.14'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.25%  '

# Row 17
$ws.Range("D17").Value = '''0.0000172'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.67%  '

# Row 18
$ws.Range("D18").Value = '''2.528.27'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.38%  '

# Row 19
$ws.Range("D19").Value = '''11.13'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.30%  '

# Row 20
$ws.Range("D20").Value = '''7.53'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.14%  '

# Row 21
$ws.Range("D21").Value = '''353.84'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.80%  '

# Row 22
$ws.Range("E22").Value = '  +5.09%  '

# Row 23
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = '''71.23'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.80%  '

# Row 24
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").Value = '''1.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.22%  '

# Row 25
$ws.Range("E25").Value = '  +1.13%  '

# Row 26
$ws.Range("E26").Value = '  -3.94%  '

# Row 27
$ws.Range("D27").Value = '''9.06'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.79%  '

# Row 28
$ws.Range("E28").Value = '  +1.71%  '

# Row 29
$ws.Range("D29").Value = '''1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.27%  '

# Row 30
$ws.Range("D30").Value = '''515.67'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.60%  '

# Row 31
$ws.Range("D31").Value = '''0.0₃0902'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.47%  '

# Row 32
$ws.Range("D32").Value = '''7.83'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.58%  '

# Row 33
$ws.Range("E33").Value = '  +2.97%  '

# Row 34
$ws.Range("E34").Value = '  +1.53%  '

# Row 35
$ws.Range("E35").Value = '  +0.02%  '

# Row 36
$ws.Range("D36").Value = '''164.17'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.44%  '

# Row 37
$ws.Range("D37").Value = '''0.120'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.70%  '

# Row 38
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").Value = '''18.45'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.98%  '

# Row 39
$ws.Range("B39").Value = 'WhiteBITCoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D39").Value = '''18.71'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.20%  '

# Row 40
$ws.Range("D40").Value = '''1.33'
$ws.Range("D40").Style = "Normal"

# Row 41
$ws.Range("E41").Value = '  +5.17%  '

# Row 42
$ws.Range("E42").Value = '  +0.06%  '

# Row 43
$ws.Range("D43").Value = '''4.87'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.32%  '

# Row 44
$ws.Range("D44").Value = '''0.328'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.58%  '

# Row 45
$ws.Range("E45").Value = '  +2.44%  '

# Row 46
$ws.Range("D46").Value = '''153.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.62%  '

# Row 47
$ws.Range("E47").Value = '  +3.16%  '

# Row 48
$ws.Range("D48").Value = '''0.0₆0262'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.89%  '

# Row 49
$ws.Range("D49").Value = '''0.522'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.74%  '

# Row 50
$ws.Range("E50").Value = '  +3.66%  '

# Row 51
$ws.Range("D51").Value = '''0.0743'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.53%  '
